# PSA-公共存储公司.xlsx -- add 2022-Q3 data
#
# Summary of the edit:
#  1. Insert a new worksheet "2022-Q3" right before the existing "2022-Q2"
#     worksheet and fill it with the Q3 fund-holding table (two new funds
#     appended compared to the Q2 sheet, several names/values refreshed).
#  2. Update the "总计" (totals) sheet: insert a new row for 2022-Q3 at the
#     top of the data and keep the rest of the historical rows, renumbering
#     the running index in column A.
#
# All other quarter sheets (2022-Q2 .. 2020-Q4) are left completely
# untouched; they simply shift one tab to the right because of the new
# sheet insertion.

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
}

function Set-NumCell($ws, $row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "General"
    $c.Value = $val
}

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" worksheet right before "2022-Q2".
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($q2)
$q3.Name = "2022-Q3"

# Header row (same headers as every other quarter sheet).
Set-TextCell $q3 1 2 "基金代码"
Set-TextCell $q3 1 3 "基金名称"
Set-TextCell $q3 1 4 "基金规模"
Set-TextCell $q3 1 5 "股票总仓位"
Set-TextCell $q3 1 6 "仓位占比"
Set-TextCell $q3 1 7 "持有市值(亿元)"
Set-TextCell $q3 1 8 "仓位排名"

$q3rows = @(
    @{ idx=0; code="000179"; name="广发美国房地产指数（QDII）人民币A"; d="1.82"; e="92.37"; f="4.64"; g="0.0844"; h=3 },
    @{ idx=1; code="000180"; name="广发美国房地产指数（QDII）美元A";   d="1.82"; e="92.37"; f="4.64"; g="0.0844"; h=3 },
    @{ idx=2; code="160140"; name="南方道琼斯美国精选REIT指数（QDII-LOF）A"; d="0.78"; e="91.13"; f="4.97"; g="0.0388"; h=3 },
    @{ idx=3; code="160141"; name="南方道琼斯美国精选REIT指数（QDII-LOF）C"; d="0.42"; e="91.13"; f="4.97"; g="0.0209"; h=3 },
    @{ idx=4; code="070031"; name="嘉实全球房地产（QDII）";             d="0.38"; e="94.39"; f="4.07"; g="0.0155"; h=3 },
    @{ idx=5; code="016278"; name="广发美国房地产指数（QDII）人民币C"; d="0.01"; e="92.37"; f="4.64"; g="0.0005"; h=3 },
    @{ idx=6; code="016279"; name="广发美国房地产指数（QDII）美元C";   d="0.01"; e="92.37"; f="4.64"; g="0.0005"; h=3 }
)

$row = 2
foreach ($r in $q3rows) {
    Set-NumCell  $q3 $row 1 $r.idx
    Set-TextCell $q3 $row 2 $r.code
    Set-TextCell $q3 $row 3 $r.name
    Set-TextCell $q3 $row 4 $r.d
    Set-TextCell $q3 $row 5 $r.e
    Set-TextCell $q3 $row 6 $r.f
    Set-TextCell $q3 $row 7 $r.g
    Set-NumCell  $q3 $row 8 $r.h
    $row++
}

# ---------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet with the new 2022-Q3 row, shifting
#    the rest of the history down by one and renumbering column A.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$totalRows = @(
    @{ idx=0; q="2022-Q3"; c=7;  d=0.24 },
    @{ idx=1; q="2022-Q2"; c=5;  d=0.28 },
    @{ idx=2; q="2022-Q1"; c=5;  d=0.33 },
    @{ idx=3; q="2021-Q4"; c=10; d=1.12 },
    @{ idx=4; q="2021-Q3"; c=8;  d=0.96 },
    @{ idx=5; q="2021-Q2"; c=7;  d=0.97 },
    @{ idx=6; q="2021-Q1"; c=8;  d=0.61 },
    @{ idx=7; q="2020-Q4"; c=7;  d=0.47 }
)

$row = 2
foreach ($r in $totalRows) {
    Set-NumCell  $total $row 1 $r.idx
    Set-TextCell $total $row 2 $r.q
    Set-NumCell  $total $row 3 $r.c
    Set-NumCell  $total $row 4 $r.d
    $row++
}

Write-Output "2022-Q3 sheet added and totals sheet updated"
